# Regenerate save_data to use K instead of Strike#: update column G (K) values
# for rows 2-11 on Sheet1.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$values = @{
    2  = 5
    3  = 2
    4  = 3
    5  = 5
    6  = 3
    7  = 4
    8  = 1
    9  = 2
    10 = 1
    11 = 3
}

foreach ($row in $values.Keys) {
    $ws.Range("G$row").Value = $values[$row]
}
